$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.34409633333334
$ws.Range("H2").Value = 94.03228900000001
$ws.Range("I2").Value = 0.1331436845910624
$ws.Range("J2").Value = 0.1436002032967755
$ws.Range("M2").Value = 13.582109
$ws.Range("N2").Value = 40.746327
$ws.Range("O2").Value = 0.1468569803870159
$ws.Range("P2").Value = 0.1573927859769708
$ws.Range("Q2").Value = 425.7189329058338
$ws.Range("R2").Value = 3831.470396152503
$ws.Range("S2").Value = 0.01955307947664469
$ws.Range("T2").Value = 0.02260163606373888
$ws.Range("G3").Value = 31.34409633333334
$ws.Range("H3").Value = 94.03228900000001
$ws.Range("I3").Value = 0.1331436845910624
$ws.Range("J3").Value = 0.1436002032967755
$ws.Range("O3").Value = 0.1979078176984568
$ws.Range("P3").Value = 0.2121061097136419
$ws.Range("Q3").Value = 573.7085478829467
$ws.Range("R3").Value = 5163.376930946521
$ws.Range("S3").Value = 0.02635017605774881
$ws.Range("T3").Value = 0.03045848047536714
$ws.Range("G4").Value = 31.34409633333334
$ws.Range("H4").Value = 94.03228900000001
$ws.Range("I4").Value = 0.1331436845910624
$ws.Range("J4").Value = 0.1436002032967755
$ws.Range("M4").Value = 20.18372733333333
$ws.Range("N4").Value = 60.551182
$ws.Range("O4").Value = 0.218237186075315
$ws.Range("P4").Value = 0.2338939465434174
$ws.Range("Q4").Value = 632.6406939017331
$ws.Range("R4").Value = 5693.766245115598
$ws.Range("S4").Value = 0.02905690306885274
$ws.Range("T4").Value = 0.03358721827351988
$ws.Range("G5").Value = 31.34409633333334
$ws.Range("H5").Value = 94.03228900000001
$ws.Range("I5").Value = 0.1331436845910624
$ws.Range("J5").Value = 0.1436002032967755
$ws.Range("M5").Value = 18.572775
$ws.Range("N5").Value = 37.14555
$ws.Range("O5").Value = 0.2008187133461717
$ws.Range("P5").Value = 0.1434838924535914
$ws.Range("Q5").Value = 582.1468487773251
$ws.Range("R5").Value = 3492.88109266395
$ws.Range("S5").Value = 0.02673774342974567
$ws.Range("T5").Value = 0.02060431612614839
$ws.Range("G6").Value = 31.34409633333334
$ws.Range("H6").Value = 94.03228900000001
$ws.Range("I6").Value = 0.1331436845910624
$ws.Range("J6").Value = 0.1436002032967755
$ws.Range("M6").Value = 21.843109
$ws.Range("N6").Value = 65.52932699999999
$ws.Range("O6").Value = 0.2361793024930407
$ws.Range("P6").Value = 0.2531232653123785
$ws.Range("Q6").Value = 684.6525127155004
$ws.Range("R6").Value = 6161.872614439503
$ws.Range("S6").Value = 0.03144578255807054
$ws.Range("T6").Value = 0.03634855235800119
$ws.Range("I7").Value = 0.1902859530282681
$ws.Range("J7").Value = 0.2052301738779898
$ws.Range("M7").Value = 13.582109
$ws.Range("N7").Value = 40.746327
$ws.Range("O7").Value = 0.1468569803870159
$ws.Range("P7").Value = 0.1573927859769708
$ws.Range("Q7").Value = 608.4279034260841
$ws.Range("R7").Value = 5475.851130834756
$ws.Range("S7").Value = 0.02794482047179701
$ws.Range("T7").Value = 0.03230174883319496
$ws.Range("I8").Value = 0.1902859530282681
$ws.Range("J8").Value = 0.2052301738779898
$ws.Range("O8").Value = 0.1979078176984568
$ws.Range("P8").Value = 0.2121061097136419
$ws.Range("S8").Value = 0.0376590777024956
$ws.Range("T8").Value = 0.04353057377711472
$ws.Range("I9").Value = 0.1902859530282681
$ws.Range("J9").Value = 0.2052301738779898
$ws.Range("M9").Value = 20.18372733333333
$ws.Range("N9").Value = 60.551182
$ws.Range("O9").Value = 0.218237186075315
$ws.Range("P9").Value = 0.2338939465434174
$ws.Range("Q9").Value = 904.155820332744
$ws.Range("R9").Value = 8137.402382994695
$ws.Range("S9").Value = 0.0415274709385488
$ws.Range("T9").Value = 0.04800209531811481
$ws.Range("I10").Value = 0.1902859530282681
$ws.Range("J10").Value = 0.2052301738779898
$ws.Range("M10").Value = 18.572775
$ws.Range("N10").Value = 37.14555
$ws.Range("O10").Value = 0.2008187133461717
$ws.Range("P10").Value = 0.1434838924535914
$ws.Range("Q10").Value = 831.9911549858999
$ws.Range("R10").Value = 4991.9469299154
$ws.Range("S10").Value = 0.03821298025498688
$ws.Range("T10").Value = 0.02944722419694135
$ws.Range("I11").Value = 0.1902859530282681
$ws.Range("J11").Value = 0.2052301738779898
$ws.Range("M11").Value = 21.843109
$ws.Range("N11").Value = 65.52932699999999
$ws.Range("O11").Value = 0.2361793024930407
$ws.Range("P11").Value = 0.2531232653123785
$ws.Range("Q11").Value = 978.489939462084
$ws.Range("R11").Value = 8806.409455158755
$ws.Range("S11").Value = 0.04494160366043987
$ws.Range("T11").Value = 0.05194853175262399
$ws.Range("G12").Value = 52.65180833333333
$ws.Range("H12").Value = 157.955425
$ws.Range("I12").Value = 0.2236547414648942
$ws.Range("J12").Value = 0.2412196000230152
$ws.Range("M12").Value = 13.582109
$ws.Range("N12").Value = 40.746327
$ws.Range("O12").Value = 0.1468569803870159
$ws.Range("P12").Value = 0.1573927859769708
$ws.Range("Q12").Value = 715.1225998304417
$ws.Range("R12").Value = 6436.103398473975
$ws.Range("S12").Value = 0.03284525998077308
$ws.Range("T12").Value = 0.03796622487987292
$ws.Range("G13").Value = 52.65180833333333
$ws.Range("H13").Value = 157.955425
$ws.Range("I13").Value = 0.2236547414648942
$ws.Range("J13").Value = 0.2412196000230152
$ws.Range("O13").Value = 0.1979078176984568
$ws.Range("P13").Value = 0.2121061097136419
$ws.Range("Q13").Value = 963.7155329376666
$ws.Range("R13").Value = 8673.439796438999
$ws.Range("S13").Value = 0.04426302180122976
$ws.Range("T13").Value = 0.05116415094756248
$ws.Range("G14").Value = 52.65180833333333
$ws.Range("H14").Value = 157.955425
$ws.Range("I14").Value = 0.2236547414648942
$ws.Range("J14").Value = 0.2412196000230152
$ws.Range("M14").Value = 20.18372733333333
$ws.Range("N14").Value = 60.551182
$ws.Range("O14").Value = 0.218237186075315
$ws.Range("P14").Value = 0.2338939465434174
$ws.Range("Q14").Value = 1062.709743006928
$ws.Range("R14").Value = 9564.38768706235
$ws.Range("S14").Value = 0.04880978142970057
$ws.Range("T14").Value = 0.05641980423300764
$ws.Range("G15").Value = 52.65180833333333
$ws.Range("H15").Value = 157.955425
$ws.Range("I15").Value = 0.2236547414648942
$ws.Range("J15").Value = 0.2412196000230152
$ws.Range("M15").Value = 18.572775
$ws.Range("N15").Value = 37.14555
$ws.Range("O15").Value = 0.2008187133461717
$ws.Range("P15").Value = 0.1434838924535914
$ws.Range("Q15").Value = 977.8901895181249
$ws.Range("R15").Value = 5867.34113710875
$ws.Range("S15").Value = 0.04491405741475073
$ws.Range("T15").Value = 0.03461112714740064
$ws.Range("G16").Value = 52.65180833333333
$ws.Range("H16").Value = 157.955425
$ws.Range("I16").Value = 0.2236547414648942
$ws.Range("J16").Value = 0.2412196000230152
$ws.Range("M16").Value = 21.843109
$ws.Range("N16").Value = 65.52932699999999
$ws.Range("O16").Value = 0.2361793024930407
$ws.Range("P16").Value = 0.2531232653123785
$ws.Range("Q16").Value = 1150.079188472108
$ws.Range("R16").Value = 10350.71269624897
$ws.Range("S16").Value = 0.05282262083844005
$ws.Range("T16").Value = 0.06105829281517149
$ws.Range("G17").Value = 51.4266815
$ws.Range("H17").Value = 102.853363
$ws.Range("I17").Value = 0.2184506386269409
$ws.Range("J17").Value = 0.1570711932425365
$ws.Range("M17").Value = 13.582109
$ws.Range("N17").Value = 40.746327
$ws.Range("O17").Value = 0.1468569803870159
$ws.Range("P17").Value = 0.1573927859769708
$ws.Range("Q17").Value = 698.4827936412836
$ws.Range("R17").Value = 4190.896761847701
$ws.Range("S17").Value = 0.03208100115236777
$ws.Range("T17").Value = 0.02472187270116998
$ws.Range("G18").Value = 51.4266815
$ws.Range("H18").Value = 102.853363
$ws.Range("I18").Value = 0.2184506386269409
$ws.Range("J18").Value = 0.1570711932425365
$ws.Range("O18").Value = 0.1979078176984568
$ws.Range("P18").Value = 0.2121061097136419
$ws.Range("Q18").Value = 941.2913504361401
$ws.Range("R18").Value = 5647.74810261684
$ws.Range("S18").Value = 0.04323308916549209
$ws.Range("T18").Value = 0.0333157597467541
$ws.Range("G19").Value = 51.4266815
$ws.Range("H19").Value = 102.853363
$ws.Range("I19").Value = 0.2184506386269409
$ws.Range("J19").Value = 0.1570711932425365
$ws.Range("M19").Value = 20.18372733333333
$ws.Range("N19").Value = 60.551182
$ws.Range("O19").Value = 0.218237186075315
$ws.Range("P19").Value = 0.2338939465434174
$ws.Range("Q19").Value = 1037.982117054178
$ws.Range("R19").Value = 6227.892702325066
$ws.Range("S19").Value = 0.0476740526702991
$ws.Range("T19").Value = 0.03673800127578063
$ws.Range("G20").Value = 51.4266815
$ws.Range("H20").Value = 102.853363
$ws.Range("I20").Value = 0.2184506386269409
$ws.Range("J20").Value = 0.1570711932425365
$ws.Range("M20").Value = 18.572775
$ws.Range("N20").Value = 37.14555
$ws.Range("O20").Value = 0.2008187133461717
$ws.Range("P20").Value = 0.1434838924535914
$ws.Range("Q20").Value = 955.1361844961625
$ws.Range("R20").Value = 3820.54473798465
$ws.Range("S20").Value = 0.04386897617871181
$ws.Range("T20").Value = 0.02253718619876938
$ws.Range("G21").Value = 51.4266815
$ws.Range("H21").Value = 102.853363
$ws.Range("I21").Value = 0.2184506386269409
$ws.Range("J21").Value = 0.1570711932425365
$ws.Range("M21").Value = 21.843109
$ws.Range("N21").Value = 65.52932699999999
$ws.Range("O21").Value = 0.2361793024930407
$ws.Range("P21").Value = 0.2531232653123785
$ws.Range("Q21").Value = 1123.318609512783
$ws.Range("R21").Value = 6739.911657076701
$ws.Range("S21").Value = 0.05159351946007021
$ws.Range("T21").Value = 0.03975837332006245
$ws.Range("G22").Value = 55.19670733333334
$ws.Range("H22").Value = 165.590122
$ws.Range("I22").Value = 0.2344649822888343
$ws.Range("J22").Value = 0.252878829559683
$ws.Range("M22").Value = 13.582109
$ws.Range("N22").Value = 40.746327
$ws.Range("O22").Value = 0.1468569803870159
$ws.Range("P22").Value = 0.1573927859769708
$ws.Range("Q22").Value = 749.6876954424328
$ws.Range("R22").Value = 6747.189258981894
$ws.Range("S22").Value = 0.03443281930543337
$ws.Range("T22").Value = 0.03980130349899406
$ws.Range("G23").Value = 55.19670733333334
$ws.Range("H23").Value = 165.590122
$ws.Range("I23").Value = 0.2344649822888343
$ws.Range("J23").Value = 0.252878829559683
$ws.Range("O23").Value = 0.1979078176984568
$ws.Range("P23").Value = 0.2121061097136419
$ws.Range("Q23").Value = 1010.296244478107
$ws.Range("R23").Value = 9092.66620030296
$ws.Range("S23").Value = 0.04640245297149051
$ws.Range("T23").Value = 0.05363714476684348
$ws.Range("G24").Value = 55.19670733333334
$ws.Range("H24").Value = 165.590122
$ws.Range("I24").Value = 0.2344649822888343
$ws.Range("J24").Value = 0.252878829559683
$ws.Range("M24").Value = 20.18372733333333
$ws.Range("N24").Value = 60.551182
$ws.Range("O24").Value = 0.218237186075315
$ws.Range("P24").Value = 0.2338939465434174
$ws.Range("Q24").Value = 1114.075290513801
$ws.Range("R24").Value = 10026.6776146242
$ws.Range("S24").Value = 0.05116897796791375
$ws.Range("T24").Value = 0.05914682744299446
$ws.Range("G25").Value = 55.19670733333334
$ws.Range("H25").Value = 165.590122
$ws.Range("I25").Value = 0.2344649822888343
$ws.Range("J25").Value = 0.252878829559683
$ws.Range("M25").Value = 18.572775
$ws.Range("N25").Value = 37.14555
$ws.Range("O25").Value = 0.2008187133461717
$ws.Range("P25").Value = 0.1434838924535914
$ws.Range("Q25").Value = 1025.15602604285
$ws.Range("R25").Value = 6150.9361562571
$ws.Range("S25").Value = 0.04708495606797664
$ws.Range("T25").Value = 0.03628403878433163
$ws.Range("G26").Value = 55.19670733333334
$ws.Range("H26").Value = 165.590122
$ws.Range("I26").Value = 0.2344649822888343
$ws.Range("J26").Value = 0.252878829559683
$ws.Range("M26").Value = 21.843109
$ws.Range("N26").Value = 65.52932699999999
$ws.Range("O26").Value = 0.2361793024930407
$ws.Range("P26").Value = 0.2531232653123785
$ws.Range("Q26").Value = 1205.667694723099
$ws.Range("R26").Value = 10851.00925250789
$ws.Range("S26").Value = 0.05537577597602002
$ws.Range("T26").Value = 0.06400951506651938

Write-Host "Done updating cells."
